$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '39.429.39'
Set-TextValue $ws.Range('E2') '  +1.69%  '
Set-TextValue $ws.Range('D3') '2.152.62'
Set-TextValue $ws.Range('E3') '  +2.33%  '
Set-TextValue $ws.Range('E4') '  +0.07%  '
Set-TextValue $ws.Range('D5') '226.95'
Set-TextValue $ws.Range('E5') '  -0.30%  '
Set-TextValue $ws.Range('D6') '0.628'
Set-TextValue $ws.Range('E6') '  +1.91%  '
Set-TextValue $ws.Range('D7') '62.79'
Set-TextValue $ws.Range('E7') '  +0.61%  '
Set-TextValue $ws.Range('E8') '  +0.07%  '
Set-TextValue $ws.Range('D9') '0.391'
Set-TextValue $ws.Range('E9') '  +0.25%  '
Set-TextValue $ws.Range('D10') '0.0845'
Set-TextValue $ws.Range('E10') '  +0.47%  '
Set-TextValue $ws.Range('D11') '0.103'
Set-TextValue $ws.Range('E11') '  +0.08%  '
Set-TextValue $ws.Range('D12') '15.89'
Set-TextValue $ws.Range('E12') '  +0.59%  '
Set-TextValue $ws.Range('D13') '2.475.86'
Set-TextValue $ws.Range('E13') '  +2.47%  '
Set-TextValue $ws.Range('D14') '21.86'
Set-TextValue $ws.Range('E14') '  -0.68%  '
Set-TextValue $ws.Range('D15') '0.803'
Set-TextValue $ws.Range('E15') '  -1.03%  '
Set-TextValue $ws.Range('D16') '5.46'
Set-TextValue $ws.Range('E16') '  -1.38%  '
Set-TextValue $ws.Range('D17') '2.156.79'
Set-TextValue $ws.Range('E17') '  +2.60%  '
Set-TextValue $ws.Range('D18') '39.417.63'
Set-TextValue $ws.Range('E18') '  +1.71%  '
Set-TextValue $ws.Range('D19') '71.79'
Set-TextValue $ws.Range('E19') '  +0.26%  '
Set-TextValue $ws.Range('D20') '6.10'
Set-TextValue $ws.Range('E20') '  -0.47%  '
Set-TextValue $ws.Range('D21') '0.0₃0841'
Set-TextValue $ws.Range('E21') '  -0.15%  '
Set-TextValue $ws.Range('D22') '227.37'
Set-TextValue $ws.Range('E22') '  -0.49%  '
Set-TextValue $ws.Range('E23') '  +0.03%  '
Set-TextValue $ws.Range('D24') '2.35'
Set-TextValue $ws.Range('E24') '  +1.82%  '
Set-TextValue $ws.Range('D25') '2.29'
Set-TextValue $ws.Range('E25') '  -2.58%  '
Set-TextValue $ws.Range('B26') 'Monero'
Set-TextValue $ws.Range('C26') 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws.Range('D26') '171.60'
Set-TextValue $ws.Range('E26') '  -0.37%  '
Set-TextValue $ws.Range('B27') 'Cosmos'
Set-TextValue $ws.Range('C27') 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue $ws.Range('D27') '9.56'
Set-TextValue $ws.Range('E27') '  -1.18%  '
Set-TextValue $ws.Range('D28') '0.138'
Set-TextValue $ws.Range('E28') '  +0.15%  '
Set-TextValue $ws.Range('D29') '19.66'
Set-TextValue $ws.Range('E29') '  +1.52%  '
Set-TextValue $ws.Range('D30') '1.40'
Set-TextValue $ws.Range('E30') '  -1.15%  '
Set-TextValue $ws.Range('D31') '2.67'
Set-TextValue $ws.Range('E31') '  +4.49%  '
Set-TextValue $ws.Range('D32') '0.122'
Set-TextValue $ws.Range('E32') '  +1.08%  '
Set-TextValue $ws.Range('D33') '4.58'
Set-TextValue $ws.Range('E33') '  +0.38%  '
Set-TextValue $ws.Range('D34') '4.67'
Set-TextValue $ws.Range('E34') '  -1.75%  '
Set-TextValue $ws.Range('D35') '6.92'
Set-TextValue $ws.Range('E35') '  -2.95%  '
Set-TextValue $ws.Range('D36') '0.0616'
Set-TextValue $ws.Range('E36') '  -0.20%  '
Set-TextValue $ws.Range('D37') '2.39'
Set-TextValue $ws.Range('E37') '  +0.07%  '
Set-TextValue $ws.Range('D38') '3.58'
Set-TextValue $ws.Range('E38') '  +0.74%  '
Set-TextValue $ws.Range('D39') '0.999'
Set-TextValue $ws.Range('E39') '  -0.18%  '
Set-TextValue $ws.Range('D40') '4.69'
Set-TextValue $ws.Range('E40') '  +12.90%  '
Set-TextValue $ws.Range('D41') '101.61'
Set-TextValue $ws.Range('E41') '  -1.08%  '
Set-TextValue $ws.Range('D42') '0.0225'
Set-TextValue $ws.Range('E42') '  -0.84%  '
Set-TextValue $ws.Range('D43') '17.60'
Set-TextValue $ws.Range('E43') '  -3.58%  '
Set-TextValue $ws.Range('D44') '1.507.70'
Set-TextValue $ws.Range('E44') '  -1.35%  '
Set-TextValue $ws.Range('D45') '1.19'
Set-TextValue $ws.Range('E45') '  -1.62%  '
Set-TextValue $ws.Range('D46') '0.0921'
Set-TextValue $ws.Range('E46') '  +0.65%  '
Set-TextValue $ws.Range('B47') 'HuobiToken'
Set-TextValue $ws.Range('C47') 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue $ws.Range('D47') '2.80'
Set-TextValue $ws.Range('E47') '  -0.11%  '
Set-TextValue $ws.Range('B48') 'FraxShare'
Set-TextValue $ws.Range('C48') 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws.Range('D48') '7.72'
Set-TextValue $ws.Range('E48') '  -0.95%  '
Set-TextValue $ws.Range('B49') 'ARBITRUM'
Set-TextValue $ws.Range('C49') 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue $ws.Range('D49') '1.09'
Set-TextValue $ws.Range('E49') '  +0.96%  '
Set-TextValue $ws.Range('B50') 'TerraClassic'
Set-TextValue $ws.Range('C50') 'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc'
Set-TextValue $ws.Range('D50') '0.000189'
Set-TextValue $ws.Range('E50') '  +36.28%  '
Set-TextValue $ws.Range('B51') 'MXToken'
Set-TextValue $ws.Range('C51') 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue $ws.Range('D51') '2.99'
Set-TextValue $ws.Range('E51') '  +0.52%  '
